$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-08-15 08:03"
$ws.Range("B2").Value = "d825236"
$ws.Range("C2").Value = "FIX: Resolve monthly target display issue by preventing useEffect from overriding local state updates"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = "Fixed monthly target display issue"
